$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.301
$ws.Range("C4").Value = 0.053
$ws.Range("D4").Value = 0.231
$ws.Range("E4").Value = 0.162
$ws.Range("G4").Value = 0.111
$ws.Range("H4").Value = 0.204
$ws.Range("J4").Value = 0.098
$ws.Range("K4").Value = 0.382
$ws.Range("L4").Value = 0.099
$ws.Range("M4").Value = 0.314
$ws.Range("N4").Value = 0.286
$ws.Range("O4").Value = 0.018
$ws.Range("P4").Value = 0.134
$ws.Range("Q4").Value = 0.574
$ws.Range("R4").Value = 0.208
$ws.Range("S4").Value = 0.456
$ws.Range("T4").Value = 0.316
$ws.Range("V4").Value = 0.297
$ws.Range("W4").Value = 0.263
$ws.Range("Y4").Value = 0.206
$ws.Range("Z4").Value = 0.467
$ws.Range("AA4").Value = 0.126
$ws.Range("AB4").Value = 0.355
$ws.Range("AE4").Value = 0.072
$ws.Range("AF4").Value = 0.728
$ws.Range("AG4").Value = 0.092
$ws.Range("AH4").Value = 0.303
$ws.Range("AI4").Value = 0.691
$ws.Range("AJ4").Value = 0.162
$ws.Range("AK4").Value = 0.402
$ws.Range("AL4").Value = 0.746
$ws.Range("AM4").Value = 0.105
$ws.Range("AN4").Value = 0.325
$ws.Range("AO4").Value = 0.722
$ws.Range("B5").Value = 0.824
$ws.Range("C5").Value = 0.145
$ws.Range("D5").Value = 0.381
$ws.Range("E5").Value = 0.706
$ws.Range("F5").Value = 0.208
$ws.Range("G5").Value = 0.456
$ws.Range("H5").Value = 0.882
$ws.Range("I5").Value = 0.104
$ws.Range("J5").Value = 0.322
$ws.Range("K5").Value = 0.735
$ws.Range("L5").Value = 0.195
$ws.Range("M5").Value = 0.441
$ws.Range("N5").Value = 0.882
$ws.Range("O5").Value = 0.104
$ws.Range("P5").Value = 0.322
$ws.Range("Q5").Value = 0.647
$ws.Range("R5").Value = 0.228
$ws.Range("S5").Value = 0.478
$ws.Range("T5").Value = 0.647
$ws.Range("U5").Value = 0.228
$ws.Range("V5").Value = 0.478
$ws.Range("W5").Value = 0.794
$ws.Range("X5").Value = 0.163
$ws.Range("Y5").Value = 0.404
$ws.Range("Z5").Value = 0.853
$ws.Range("AA5").Value = 0.125
$ws.Range("AB5").Value = 0.354
$ws.Range("AC5").Value = 0.824
$ws.Range("AD5").Value = 0.145
$ws.Range("AE5").Value = 0.381
$ws.Range("AF5").Value = 0.971
$ws.Range("AH5").Value = 0.169
$ws.Range("AI5").Value = 0.794
$ws.Range("AJ5").Value = 0.163
$ws.Range("AK5").Value = 0.404
$ws.Range("AL5").Value = 0.941
$ws.Range("AM5").Value = 0.055
$ws.Range("AN5").Value = 0.235
$ws.Range("AO5").Value = 0.902
$ws.Range("B6").Value = 0.441
$ws.Range("E6").Value = 0.264
$ws.Range("H6").Value = 0.331
$ws.Range("K6").Value = 0.503
$ws.Range("N6").Value = 0.432
$ws.Range("Q6").Value = 0.608
$ws.Range("T6").Value = 0.425
$ws.Range("W6").Value = 0.395
$ws.Range("Z6").Value = 0.604
$ws.Range("AF6").Value = 0.832
$ws.Range("AI6").Value = 0.739
$ws.Range("AL6").Value = 0.832
$ws.Range("AO6").Value = 0.801
$ws.Range("B7").Value = 0.611
$ws.Range("E7").Value = 0.422
$ws.Range("H7").Value = 0.53
$ws.Range("K7").Value = 0.62
$ws.Range("N7").Value = 0.623
$ws.Range("Q7").Value = 0.631
$ws.Range("T7").Value = 0.535
$ws.Range("W7").Value = 0.566
$ws.Range("Z7").Value = 0.732
$ws.Range("AC7").Value = 0.408
$ws.Range("AF7").Value = 0.91
$ws.Range("AI7").Value = 0.771
$ws.Range("AL7").Value = 0.894
$ws.Range("AO7").Value = 0.858
$ws.Range("B8").Value = 0.762
$ws.Range("C8").Value = 0.146
$ws.Range("D8").Value = 0.383
$ws.Range("E8").Value = 0.594
$ws.Range("F8").Value = 0.184
$ws.Range("G8").Value = 0.429
$ws.Range("H8").Value = 0.768
$ws.Range("I8").Value = 0.122
$ws.Range("J8").Value = 0.349
$ws.Range("K8").Value = 0.655
$ws.Range("L8").Value = 0.183
$ws.Range("M8").Value = 0.428
$ws.Range("N8").Value = 0.793
$ws.Range("O8").Value = 0.115
$ws.Range("P8").Value = 0.339
$ws.Range("Q8").Value = 0.614
$ws.Range("R8").Value = 0.216
$ws.Range("S8").Value = 0.465
$ws.Range("T8").Value = 0.56
$ws.Range("U8").Value = 0.197
$ws.Range("V8").Value = 0.444
$ws.Range("W8").Value = 0.71
$ws.Range("X8").Value = 0.157
$ws.Range("Y8").Value = 0.396
$ws.Range("Z8").Value = 0.791
$ws.Range("AA8").Value = 0.13
$ws.Range("AB8").Value = 0.361
$ws.Range("AC8").Value = 0.702
$ws.Range("AD8").Value = 0.154
$ws.Range("AE8").Value = 0.392
$ws.Range("AF8").Value = 0.891
$ws.Range("AG8").Value = 0.049
$ws.Range("AH8").Value = 0.221
$ws.Range("AI8").Value = 0.783
$ws.Range("AJ8").Value = 0.163
$ws.Range("AK8").Value = 0.404
$ws.Range("AL8").Value = 0.909
$ws.Range("AM8").Value = 0.062
$ws.Range("AN8").Value = 0.25
$ws.Range("AO8").Value = 0.861
$ws.Range("B9").Value = 0.676
$ws.Range("C9").Value = 0.219
$ws.Range("D9").Value = 0.468
$ws.Range("E9").Value = 0.471
$ws.Range("F9").Value = 0.249
$ws.Range("G9").Value = 0.499
$ws.Range("H9").Value = 0.647
$ws.Range("I9").Value = 0.228
$ws.Range("J9").Value = 0.478
$ws.Range("K9").Value = 0.559
$ws.Range("L9").Value = 0.247
$ws.Range("M9").Value = 0.497
$ws.Range("N9").Value = 0.676
$ws.Range("O9").Value = 0.219
$ws.Range("P9").Value = 0.468
$ws.Range("Q9").Value = 0.559
$ws.Range("R9").Value = 0.247
$ws.Range("S9").Value = 0.497
$ws.Range("T9").Value = 0.441
$ws.Range("U9").Value = 0.247
$ws.Range("V9").Value = 0.497
$ws.Range("W9").Value = 0.588
$ws.Range("X9").Value = 0.242
$ws.Range("Y9").Value = 0.492
$ws.Range("Z9").Value = 0.706
$ws.Range("AA9").Value = 0.208
$ws.Range("AB9").Value = 0.456
$ws.Range("AC9").Value = 0.588
$ws.Range("AD9").Value = 0.242
$ws.Range("AE9").Value = 0.492
$ws.Range("AF9").Value = 0.765
$ws.Range("AG9").Value = 0.18
$ws.Range("AH9").Value = 0.424
$ws.Range("AI9").Value = 0.765
$ws.Range("AJ9").Value = 0.18
$ws.Range("AK9").Value = 0.424
$ws.Range("AL9").Value = 0.853
$ws.Range("AM9").Value = 0.125
$ws.Range("AN9").Value = 0.354
$ws.Range("AO9").Value = 0.794
$ws.Range("B10").Value = 0.824
$ws.Range("C10").Value = 0.145
$ws.Range("D10").Value = 0.381
$ws.Range("E10").Value = 0.618
$ws.Range("F10").Value = 0.236
$ws.Range("G10").Value = 0.486
$ws.Range("H10").Value = 0.794
$ws.Range("I10").Value = 0.163
$ws.Range("J10").Value = 0.404
$ws.Range("K10").Value = 0.735
$ws.Range("L10").Value = 0.195
$ws.Range("M10").Value = 0.441
$ws.Range("N10").Value = 0.853
$ws.Range("O10").Value = 0.125
$ws.Range("P10").Value = 0.354
$ws.Range("Q10").Value = 0.647
$ws.Range("R10").Value = 0.228
$ws.Range("S10").Value = 0.478
$ws.Range("T10").Value = 0.647
$ws.Range("U10").Value = 0.228
$ws.Range("V10").Value = 0.478
$ws.Range("W10").Value = 0.794
$ws.Range("X10").Value = 0.163
$ws.Range("Y10").Value = 0.404
$ws.Range("Z10").Value = 0.853
$ws.Range("AA10").Value = 0.125
$ws.Range("AB10").Value = 0.354
$ws.Range("AC10").Value = 0.706
$ws.Range("AD10").Value = 0.208
$ws.Range("AE10").Value = 0.456
$ws.Range("AF10").Value = 0.971
$ws.Range("AH10").Value = 0.169
$ws.Range("AI10").Value = 0.794
$ws.Range("AJ10").Value = 0.163
$ws.Range("AK10").Value = 0.404
$ws.Range("AL10").Value = 0.941
$ws.Range("AM10").Value = 0.055
$ws.Range("AN10").Value = 0.235
$ws.Range("AO10").Value = 0.902
$ws.Range("B11").Value = 0.824
$ws.Range("C11").Value = 0.145
$ws.Range("D11").Value = 0.381
$ws.Range("E11").Value = 0.706
$ws.Range("F11").Value = 0.208
$ws.Range("G11").Value = 0.456
$ws.Range("H11").Value = 0.882
$ws.Range("I11").Value = 0.104
$ws.Range("J11").Value = 0.322
$ws.Range("K11").Value = 0.735
$ws.Range("L11").Value = 0.195
$ws.Range("M11").Value = 0.441
$ws.Range("N11").Value = 0.882
$ws.Range("O11").Value = 0.104
$ws.Range("P11").Value = 0.322
$ws.Range("Q11").Value = 0.647
$ws.Range("R11").Value = 0.228
$ws.Range("S11").Value = 0.478
$ws.Range("T11").Value = 0.647
$ws.Range("U11").Value = 0.228
$ws.Range("V11").Value = 0.478
$ws.Range("W11").Value = 0.794
$ws.Range("X11").Value = 0.163
$ws.Range("Y11").Value = 0.404
$ws.Range("Z11").Value = 0.853
$ws.Range("AA11").Value = 0.125
$ws.Range("AB11").Value = 0.354
$ws.Range("AC11").Value = 0.765
$ws.Range("AD11").Value = 0.18
$ws.Range("AE11").Value = 0.424
$ws.Range("AF11").Value = 0.971
$ws.Range("AH11").Value = 0.169
$ws.Range("AI11").Value = 0.794
$ws.Range("AJ11").Value = 0.163
$ws.Range("AK11").Value = 0.404
$ws.Range("AL11").Value = 0.941
$ws.Range("AM11").Value = 0.055
$ws.Range("AN11").Value = 0.235
$ws.Range("AO11").Value = 0.902
$ws.Range("B12").Value = 1.25
$ws.Range("C12").Value = 0.33
$ws.Range("D12").Value = 0.575
$ws.Range("E12").Value = 1.667
$ws.Range("F12").Value = 1.139
$ws.Range("G12").Value = 1.067
$ws.Range("H12").Value = 1.6
$ws.Range("I12").Value = 1.373
$ws.Range("J12").Value = 1.172
$ws.Range("K12").Value = 1.4
$ws.Range("L12").Value = 0.56
$ws.Range("M12").Value = 0.748
$ws.Range("N12").Value = 1.367
$ws.Range("O12").Value = 0.566
$ws.Range("P12").Value = 0.752
$ws.Range("Z12").Value = 1.241
$ws.Range("AA12").Value = 0.321
$ws.Range("AB12").Value = 0.567
$ws.Range("AC12").Value = 1.821
$ws.Range("AD12").Value = 2.504
$ws.Range("AE12").Value = 1.582
$ws.Range("AF12").Value = 1.242
$ws.Range("AG12").Value = 0.244
$ws.Range("AH12").Value = 0.494
$ws.Range("AI12").Value = 1.037
$ws.Range("AJ12").Value = 0.036
$ws.Range("AK12").Value = 0.189
$ws.Range("AL12").Value = 1.094
$ws.Range("AM12").Value = 0.085
$ws.Range("AN12").Value = 0.291
$ws.Range("AO12").Value = 1.124
$ws.Range("B13").Value = 3.441
$ws.Range("C13").Value = 1.423
$ws.Range("D13").Value = 1.193
$ws.Range("E13").Value = 4.571
$ws.Range("F13").Value = 0.459
$ws.Range("G13").Value = 0.678
$ws.Range("H13").Value = 4.594
$ws.Range("I13").Value = 0.679
$ws.Range("J13").Value = 0.824
$ws.Range("K13").Value = 2.265
$ws.Range("L13").Value = 0.606
$ws.Range("M13").Value = 0.779
$ws.Range("N13").Value = 3.235
$ws.Range("O13").Value = 0.768
$ws.Range("P13").Value = 0.876
$ws.Range("Z13").Value = 2.515
$ws.Range("AA13").Value = 2.916
$ws.Range("AB13").Value = 1.708
$ws.Range("AC13").Value = 6.353
$ws.Range("AD13").Value = 2.228
$ws.Range("AE13").Value = 1.493
$ws.Range("AF13").Value = 1.588
$ws.Range("AG13").Value = 0.595
$ws.Range("AH13").Value = 0.771
$ws.Range("AI13").Value = 1.206
$ws.Range("AJ13").Value = 0.163
$ws.Range("AK13").Value = 0.404
$ws.Range("AL13").Value = 1.5
$ws.Range("AM13").Value = 0.721
$ws.Range("AN13").Value = 0.849
$ws.Range("AO13").Value = 1.431
